$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) cells whose new values look numeric: force Text format first
# so Excel stores the exact original string (matching the source data feed),
# instead of silently converting to a floating point number.
$numericLookingPriceCells = @("D5","D6","D8","D11","D14","D18","D19","D20","D21","D23","D24","D26","D27","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D46","D48","D49","D50","D51")
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Updated Price values (column D)
$ws.Range("D2").Value = "63.292.02"
$ws.Range("D3").Value = "2.680.44"
$ws.Range("D5").Value = "552.73"
$ws.Range("D6").Value = "157.71"
$ws.Range("D8").Value = "0.583"
$ws.Range("D11").Value = "0.368"
$ws.Range("D13").Value = "3.152.30"
$ws.Range("D14").Value = "26.25"
$ws.Range("D15").Value = "63.165.16"
$ws.Range("D17").Value = "2.677.73"
$ws.Range("D18").Value = "11.95"
$ws.Range("D19").Value = "4.55"
$ws.Range("D20").Value = "343.42"
$ws.Range("D21").Value = "6.30"
$ws.Range("D23").Value = "0.506"
$ws.Range("D24").Value = "63.66"
$ws.Range("D26").Value = "0.997"
$ws.Range("D27").Value = "8.14"
$ws.Range("D28").Value = "0.0₃0851"
$ws.Range("D30").Value = "1.34"
$ws.Range("D32").Value = "165.77"
$ws.Range("D33").Value = "0.999"
$ws.Range("D34").Value = "4.78"
$ws.Range("D35").Value = "19.49"
$ws.Range("D36").Value = "1.43"
$ws.Range("D37").Value = "1.78"
$ws.Range("D38").Value = "339.53"
$ws.Range("D39").Value = "0.942"
$ws.Range("D40").Value = "6.07"
$ws.Range("D41").Value = "38.08"
$ws.Range("D43").Value = "20.25"
$ws.Range("D44").Value = "20.70"
$ws.Range("D46").Value = "0.0560"
$ws.Range("D48").Value = "11.05"
$ws.Range("D49").Value = "129.26"
$ws.Range("D50").Value = "0.0969"
$ws.Range("D51").Value = "0.0241"

# --- Updated Volume(1h) values (column E)
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -3.82%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("E9").Value = "  -4.02%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("E12").Value = "  -9.32%  "
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("E16").Value = "  -4.28%  "
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  -5.55%  "
$ws.Range("E20").Value = "  -4.45%  "
$ws.Range("E21").Value = "  -4.92%  "
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("E28").Value = "  -5.77%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("E39").Value = "  -6.02%  "
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("E42").Value = "  -6.60%  "
$ws.Range("E43").Value = "  -5.91%  "
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -6.09%  "
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("E51").Value = "  -4.69%  "

